$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Swap Id (A), Ost (Q) and Nord (R) values between row 2 and row 4
$a2 = $ws.Range("A2").Value2
$q2 = $ws.Range("Q2").Value2
$r2 = $ws.Range("R2").Value2

$a4 = $ws.Range("A4").Value2
$q4 = $ws.Range("Q4").Value2
$r4 = $ws.Range("R4").Value2

$ws.Range("A2").Value = $a4
$ws.Range("Q2").Value = $q4
$ws.Range("R2").Value = $r4

$ws.Range("A4").Value = $a2
$ws.Range("Q4").Value = $q2
$ws.Range("R4").Value = $r2
